$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be created in this exact order (matches the ---
# --- author's commit so the sharedStrings table lines up index-for-index) ---
$ws.Range("C4").Value = " We\'ve been away a while,\nand we just got back."
$ws.Range("C5").Value = " While we were gone, I hear\n[CS:N]Wigglytuff[CR]\'s Guild went on an expedition?"
$ws.Range("C6").Value = " Tch! That\'s disappointing.[K]\nIf we would have known about it, we would\nhave asked to tag along."

$ws.Range("A4").Value = "SCRIPT/G01P03A/um1103.ssb"

$ws.Range("D4").Value = " Нас не было некоторое время\nи мы только что вернулись."
$ws.Range("D5").Value = " Пока нас не было, Гильдия\n[CS:N]Виглитаффа[CR] отправлялась в экспедицию?"
$ws.Range("D6").Value = " Эх! Какая жалость.[K] Если бы мы\nо ней знали, мы бы напросились к вам."

$ws.Range("E4").Value = " Îàò îå áúìï îåëïóïñïå âñåíÿ\né îú óïìûëï œóï âåñîôìéòû."
$ws.Range("E5").Value = " Ðïëà îàò îå áúìï, Ãéìûäéÿ\n[CS:N]Âéãìéóàõõà[CR] ïóðñàâìÿìàòû â üëòðåäéøéý?"
$ws.Range("E6").Value = " Üö! Ëàëàÿ çàìïòóû.[K] Åòìé áú íú\nï îåê èîàìé, íú áú îàðñïòéìéòû ë âàí."

# --- Line numbers (plain numeric cells, column B) ---
$ws.Range("B4").Value = 183
$ws.Range("B5").Value = 186
$ws.Range("B6").Value = 189

# --- Row heights to match the wrapped-text content ---
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 31.8
$ws.Rows.Item(6).RowHeight = 31.8

# --- Close out the previous entry's row (row 3) with a bottom border ---
$ws.Range("A3:E3").Borders.Item(9).LineStyle = 1
$ws.Range("A3:E3").Borders.Item(9).Weight = 2

# --- Match the author's final selection ---
[void]$ws.Range("C1").Select()
